# Add a new bulleted explanation paragraph right after the last table
# (the toString()/toFixed() union-type example), before the existing
# trailing empty paragraphs at the end of the document.

$d = $word.ActiveDocument

# Locate the last table in the document (the toString()/toFixed() example)
# and build a collapsed Range immediately after it.
$tableCount = $d.Tables.Count
$lastTable = $d.Tables($tableCount)
$tblRange = $lastTable.Range
$insertionPoint = $d.Range($tblRange.End, $tblRange.End)

# WordprocessingML fragment: an empty paragraph, the new list paragraph
# (ListParagraph style, numId 7) with mixed bold/plain runs, then two more
# empty paragraphs. InsertXML merges the *last* paragraph of the fragment
# into the paragraph already sitting at the insertion point, so ending the
# fragment with an empty <w:p/> leaves that existing (empty) paragraph
# untouched while everything before it becomes brand-new paragraphs.
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Here, toString() exists for both type </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Boolean </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>string</w:t></w:r><w:r><w:t xml:space="preserve">, but </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>toFixed()</w:t></w:r><w:r><w:t xml:space="preserve"> doesn&#8217;t exist in both</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$null = $insertionPoint.InsertXML($xml)
